$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Patients sheet: add two new patient rows (17, 18)
# ---------------------------------------------------------------------
$patients = $wb.Worksheets("Patients")

$patients.Range("A16").Copy()
$patients.Range("A17").PasteSpecial(-4122)
$patients.Cells.Item(17, 1).Value = 15
$patients.Cells.Item(17, 2).Value = 16
$patients.Cells.Item(17, 3).Value = "Random patient"
$patients.Cells.Item(17, 4).Value = 28
$patients.Cells.Item(17, 5).Value = "Male"
$patients.Cells.Item(17, 6).Value = 89
$patients.Cells.Item(17, 7).Value = "obese, unhealthy"
$patients.Cells.Item(17, 8).Value = "2023-01-18 16:13:51.803629"

$patients.Range("A16").Copy()
$patients.Range("A18").PasteSpecial(-4122)
$patients.Cells.Item(18, 1).Value = 16
$patients.Cells.Item(18, 2).Value = 17
$patients.Cells.Item(18, 3).Value = 'anup`'
$patients.Cells.Item(18, 4).Value = 110
$patients.Cells.Item(18, 5).Value = "Male"
$patients.Cells.Item(18, 6).Value = 95
$patients.Cells.Item(18, 7).Value = "diabitics"
$patients.Cells.Item(18, 8).Value = "2023-01-20 11:46:43.932040"

# ---------------------------------------------------------------------
# Drugs sheet: add one new (mostly empty) drug row (13)
# ---------------------------------------------------------------------
$drugs = $wb.Worksheets("Drugs")

$drugs.Range("A12").Copy()
$drugs.Range("A13").PasteSpecial(-4122)
$drugs.Cells.Item(13, 1).Value = 11
$drugs.Cells.Item(13, 2).Value = 12
$drugs.Cells.Item(13, 7).Value = "2023-01-18 16:14:00.893199"

# ---------------------------------------------------------------------
# Record sheet: add four new record rows (17, 18, 19, 20)
# ---------------------------------------------------------------------
$record = $wb.Worksheets("Record")

$record.Range("A16").Copy()
$record.Range("A17").PasteSpecial(-4122)
$record.Cells.Item(17, 1).Value = 15
$record.Cells.Item(17, 2).Value = "Amogha"
$record.Cells.Item(17, 3).Value = "bleeding,vomit"
$record.Cells.Item(17, 4).Value = 1
$record.Cells.Item(17, 5).Value = 2
$record.Cells.Item(17, 9).Value = "2023-01-18 14:18:43.047424"

$record.Range("A16").Copy()
$record.Range("A18").PasteSpecial(-4122)
$record.Cells.Item(18, 1).Value = 16
$record.Cells.Item(18, 2).Value = "Dheeraj"
$record.Cells.Item(18, 3).Value = "bleeding,vomit,motions"
$record.Cells.Item(18, 4).Value = 9
$record.Cells.Item(18, 5).Value = 1
$record.Cells.Item(18, 6).Value = "treatment"
$record.Cells.Item(18, 7).Value = "side effects"
$record.Cells.Item(18, 8).Value = "drugs used"
$record.Cells.Item(18, 9).Value = "2023-01-18 16:08:08.807079"

$record.Range("A16").Copy()
$record.Range("A19").PasteSpecial(-4122)
$record.Cells.Item(19, 1).Value = 17
$record.Cells.Item(19, 2).Value = "Uday"
# Symptom text is identical to the existing C11 entry (trailing newline); copy
# the value across instead of retyping it so the engine reuses the shared
# string and doesn't trigger an auto row-height recalculation.
$record.Range("C11").Copy()
$record.Range("C19").PasteSpecial(-4163)
$record.Cells.Item(19, 4).Value = 7
$record.Cells.Item(19, 5).Value = 3
$record.Cells.Item(19, 6).Value = "treatment"
$record.Cells.Item(19, 7).Value = "side effects"
$record.Cells.Item(19, 8).Value = "drugs used"
$record.Cells.Item(19, 9).Value = "2023-01-18 16:08:23.286546"

$record.Range("A16").Copy()
$record.Range("A20").PasteSpecial(-4122)
$record.Cells.Item(20, 1).Value = 18
$record.Cells.Item(20, 2).Value = "Amogha"
$record.Cells.Item(20, 4).Value = 1
$record.Cells.Item(20, 5).Value = 1
$record.Cells.Item(20, 9).Value = "2023-01-18 16:13:56.863636"
